# Insert a new weekly price record as row 190 in the "Haba" sheet.
# This shifts all existing rows from 190-229 down to 191-230 (unchanged),
# and grows the used range to A1:R230.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 190, pushing rows 190-229 down to 191-230.
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with the new record's data.
$ws.Cells.Item(190, 1).Value2 = 3
$ws.Cells.Item(190, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(190, 3).Value2 = "Coquimbo"
$ws.Cells.Item(190, 4).Value2 = 44889
$ws.Cells.Item(190, 5).Value2 = 5
$ws.Cells.Item(190, 6).Value2 = 100112026
$ws.Cells.Item(190, 7).Value2 = "Haba"
$ws.Cells.Item(190, 8).Value2 = "Sin especificar"
$ws.Cells.Item(190, 9).Value2 = "Primera"
$ws.Cells.Item(190, 10).Value2 = 85
$ws.Cells.Item(190, 11).Value2 = 8500
$ws.Cells.Item(190, 12).Value2 = 9000
$ws.Cells.Item(190, 13).Value2 = 8765
$ws.Cells.Item(190, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(190, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(190, 16).Value2 = 351
$ws.Cells.Item(190, 17).Value2 = 25
$ws.Cells.Item(190, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the other
# rows in column D (style index 2 in the original workbook).
$ws.Cells.Item(190, 4).NumberFormat = $ws.Cells.Item(191, 4).NumberFormat
